# Rename worksheets and update their header-label cells to match
# the new naming convention described in the diff:
#   cw_summary   -> commission_summary_sheet  (A1 text: commissions_summary_list -> commission_summary_list)
#   cw_signups   -> signup_sheet              (A1 text: signup_commissions_table -> signup_commission_table)
#   cw_purchases -> purchase_sheet            (A1 text: purchase_commissions_table -> purchase_commission_table)

$wb = $excel.ActiveWorkbook

$summarySheet  = $wb.Worksheets.Item("cw_summary")
$signupSheet   = $wb.Worksheets.Item("cw_signups")
$purchaseSheet = $wb.Worksheets.Item("cw_purchases")

$summarySheet.Range("A1").Value  = "commission_summary_list"
$signupSheet.Range("A1").Value   = "signup_commission_table"
$purchaseSheet.Range("A1").Value = "purchase_commission_table"

$summarySheet.Name  = "commission_summary_sheet"
$signupSheet.Name   = "signup_sheet"
$purchaseSheet.Name = "purchase_sheet"
